$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three oldest "Estado de Cuenta" data rows (periods 2505, 2506,
# 2507), keeping only the most recent one (2508). That row shifts up from
# row 19 to row 16, and the summary rows below it (24/25) shift up to 21/22.
$ws.Rows("16:18").Delete()

# Update "Antigua BD" values: the overdue amount and the period count.
$ws.Range("E11").Value = 57200
$ws.Range("F13").Value = 1
